# Applies the "Add files via upload" commit to the wales_cymru-premier
# 2023-2024 workbook:
#   1) Eight match rows (spreadsheet rows 15-17 and 19-23) get their
#      match-detail columns (F:V) re-shuffled into a different row order
#      (the row's Indice/pais/torneio/temporada/data_partida in A:E stay
#      put; only the F:V "payload" for each fixture moves).
#   2) Five brand-new fixture rows are appended at the end (rows 48-52).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: re-shuffle the F:V payload of rows 15-17 & 19-23.
# Read every source row first (before anything is overwritten), then
# write them all back out in their new positions.
# ---------------------------------------------------------------------
$row15 = $ws.Range("F15:V15").Value2
$row16 = $ws.Range("F16:V16").Value2
$row17 = $ws.Range("F17:V17").Value2
$row19 = $ws.Range("F19:V19").Value2
$row20 = $ws.Range("F20:V20").Value2
$row21 = $ws.Range("F21:V21").Value2
$row22 = $ws.Range("F22:V22").Value2
$row23 = $ws.Range("F23:V23").Value2

$ws.Range("F15:V15").Value2 = $row17
$ws.Range("F16:V16").Value2 = $row15
$ws.Range("F17:V17").Value2 = $row16

$ws.Range("F19:V19").Value2 = $row22
$ws.Range("F20:V20").Value2 = $row23
$ws.Range("F21:V21").Value2 = $row20
$ws.Range("F22:V22").Value2 = $row19
$ws.Range("F23:V23").Value2 = $row21

# ---------------------------------------------------------------------
# Step 2: append five new fixture rows (48-52), matching the existing
# column layout / types / formatting used by the rest of the sheet.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=48; Idx=47; F="Caernarfon";           G=1; H="TNS";                 I=3;
       J=16.08; K="25/09/2023 16:42"; L=14.33;              M="26/09/2023 20:40";
       N=10.4;  O="25/09/2023 16:42"; P=9.050000000000001;  Q="26/09/2023 20:40";
       R=1.08;  S="25/09/2023 16:42"; T=1.14;                U="26/09/2023 20:39";
       V="https://www.betexplorer.com/football/wales/cymru-premier/caernarfon-tns/OCrsm446/" },

    @{ Row=49; Idx=48; F="Cardiff Metropolitan"; G=1; H="Haverfordwest";       I=1;
       J=2.08;  K="25/09/2023 16:42"; L=2.05;               M="26/09/2023 18:50";
       N=3.48;  O="25/09/2023 16:42"; P=3.55;               Q="26/09/2023 19:55";
       R=3.09;  S="25/09/2023 16:42"; T=3.48;               U="26/09/2023 18:50";
       V="https://www.betexplorer.com/football/wales/cymru-premier/cardiff-metropolitan-university-haverfordwest/6yYonOJC/" },

    @{ Row=50; Idx=49; F="Connahs Q.";           G=2; H="Colwyn Bay";         I=1;
       J=1.25;  K="25/09/2023 16:42"; L=1.17;               M="26/09/2023 18:51";
       N=5.62;  O="25/09/2023 16:42"; P=7.34;               Q="26/09/2023 20:22";
       R=8.17;  S="25/09/2023 16:42"; T=15.32;              U="26/09/2023 20:22";
       V="https://www.betexplorer.com/football/wales/cymru-premier/connahs-q-colwyn-bay/pSWkorZI/" },

    @{ Row=51; Idx=50; F="Penybont";             G=1; H="Barry";             I=0;
       J=1.39;  K="25/09/2023 16:42"; L=1.32;               M="26/09/2023 11:21";
       N=4.51;  O="25/09/2023 16:42"; P=5.32;               Q="26/09/2023 20:41";
       R=6.32;  S="25/09/2023 16:42"; T=9.220000000000001;  U="26/09/2023 20:41";
       V="https://www.betexplorer.com/football/wales/cymru-premier/penybont-barry-town/YsQX9RBJ/" },

    @{ Row=52; Idx=51; F="Pontypridd";           G=0; H="Newtown";           I=1;
       J=2.71;  K="25/09/2023 16:42"; L=3.15;               M="26/09/2023 18:24";
       N=3.15;  O="25/09/2023 16:42"; P=3.26;               Q="26/09/2023 18:47";
       R=2.41;  S="25/09/2023 16:42"; T=2.28;               U="26/09/2023 18:24";
       V="https://www.betexplorer.com/football/wales/cymru-premier/pontypridd-united-newtown/CYOy9oRP/" }
)

# All five new fixtures share the same kickoff date/time: 2023-09-26 20:45.
$matchDateSerial = 45195.86458333334

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # A: Indice (bold, centred, thin-bordered - matches the rest of column A)
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value2 = $r.Idx
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160
    $cellA.Borders.LineStyle = 1

    # B:D: pais / torneio / temporada
    $ws.Cells.Item($rowNum, 2).Value2 = "wales"
    $ws.Cells.Item($rowNum, 3).Value2 = "cymru-premier"
    $ws.Cells.Item($rowNum, 4).Value2 = "2023-2024"

    # E: data_partida (date-time serial, formatted like the rest of col E)
    $cellE = $ws.Cells.Item($rowNum, 5)
    $cellE.Value2 = $matchDateSerial
    $cellE.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # F:V: match payload
    $ws.Cells.Item($rowNum, 6).Value2 = $r.F
    $ws.Cells.Item($rowNum, 7).Value2 = $r.G
    $ws.Cells.Item($rowNum, 8).Value2 = $r.H
    $ws.Cells.Item($rowNum, 9).Value2 = $r.I
    $ws.Cells.Item($rowNum, 10).Value2 = $r.J
    $ws.Cells.Item($rowNum, 11).Value2 = $r.K
    $ws.Cells.Item($rowNum, 12).Value2 = $r.L
    $ws.Cells.Item($rowNum, 13).Value2 = $r.M
    $ws.Cells.Item($rowNum, 14).Value2 = $r.N
    $ws.Cells.Item($rowNum, 15).Value2 = $r.O
    $ws.Cells.Item($rowNum, 16).Value2 = $r.P
    $ws.Cells.Item($rowNum, 17).Value2 = $r.Q
    $ws.Cells.Item($rowNum, 18).Value2 = $r.R
    $ws.Cells.Item($rowNum, 19).Value2 = $r.S
    $ws.Cells.Item($rowNum, 20).Value2 = $r.T
    $ws.Cells.Item($rowNum, 21).Value2 = $r.U
    $ws.Cells.Item($rowNum, 22).Value2 = $r.V
}
